# Adicionando mensagens de status de criação de conta ao terminal
# Appends the new "usuarios" rows (ids 2-7) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - id 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "usuarioteste@gmail.com"
$ws.Range("C3").Value = "usuarioteste"
$ws.Range("D3").Value = "'10/12/2004"
$ws.Range("E3").Value = "Prefiro não informar"
$ws.Range("F3").Value = "Alternativa, Experimental"

# Row 4 - id 3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "aaaaaaaa@gmail.com"
$ws.Range("C4").Value = "aaaaaaa"
$ws.Range("D4").Value = "18/11/2004"
$ws.Range("E4").Value = "Feminino"
$ws.Range("F4").Value = "Pop, Alternativa"

# Row 5 - id 4
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "fgh@gmail.com"
$ws.Range("C5").Value = "qwsdefrf"
$ws.Range("D5").Value = "'11/11/2011"
$ws.Range("E5").Value = "Feminino"
$ws.Range("F5").Value = "Grunge"

# Row 6 - id 5
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "asdf@gmail.com"
$ws.Range("C6").Value = "ad"
$ws.Range("D6").Value = "'11/11/2011"
$ws.Range("E6").Value = "Feminino"
$ws.Range("F6").Value = "Pop"

# Row 7 - id 6
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "aaaaaa@gmail.com"
$ws.Range("C7").Value = "a"
$ws.Range("D7").Value = "'11/11/2004"
$ws.Range("E7").Value = "Feminino"
$ws.Range("F7").Value = "Rock"

# Row 8 - id 7
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "ghjk@gmail.com"
$ws.Range("C8").Value = "sedrfgt"
$ws.Range("D8").Value = "'11/11/2004"
$ws.Range("E8").Value = "Feminino"
$ws.Range("F8").Value = "Hip-Hop"

Write-Host "Usuarios adicionados com sucesso: ids 2-7 (linhas 3-8)."
